$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.486.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3658"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3344"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07445"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.994"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.925"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.574.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06752"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.399"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.482.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.387"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.624"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.019"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.752.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.177"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.811"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08281"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02437"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2265"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06488"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.415"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6079"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.756"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.051"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.56%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.221"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07222"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.43%  "
